$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "New York University, Tandon School of Engineering"
#           -> "New York University"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("New York University, Tandon School of Engineering", $true, $false, $false, $false, $false,
                         $true, 1, $false, "New York University", 2)

# ---------------------------------------------------------------------
# Change 2: merge "Incorporated a 4-D consensus cross attention module"
#           + " in " into a single run (text unchanged, runs merge)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Incorporated a 4-D consensus cross attention module in ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Incorporated a 4-D consensus cross attention module in ", 2)

# ---------------------------------------------------------------------
# Change 3: merge " to align query and support " + "features" +
#           ", increasing accuracy by 5.4%" into a single run
# ---------------------------------------------------------------------
$d.Content.Find.Execute(" to align query and support features, increasing accuracy by 5.4%", $true, $false, $false, $false, $false,
                         $true, 1, $false, " to align query and support features, increasing accuracy by 5.4%", 2)

# ---------------------------------------------------------------------
# Change 4: merge "Built" + " a scalable " into a single run
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Built a scalable ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Built a scalable ", 2)

# ---------------------------------------------------------------------
# Change 5a: merge "PyTorch" + " " (both bold/italic) into "PyTorch "
#            The word "PyTorch" also appears earlier in this same
#            paragraph-group (bold/italic, but followed by plain text),
#            so we narrow to a precise Range before searching/replacing
#            to avoid touching that other occurrence / crossing the
#            bold->plain formatting boundary.
# ---------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("Built a scalable ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Collapse(0)
$anchor.Find.Execute("PyTorch", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.MoveEnd(1, 1)
$narrow = $d.Range($anchor.Start, $anchor.End)
$narrow.Find.Execute("PyTorch ", $false, $false, $false, $false, $false, $true, 0, $false, "PyTorch ", 2)

# ---------------------------------------------------------------------
# Change 5b: merge "codebase for " + "running " +
#            "few-shot segmentation research" + " experiments " into a
#            single run
# ---------------------------------------------------------------------
$d.Content.Find.Execute("codebase for running few-shot segmentation research experiments ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "codebase for running few-shot segmentation research experiments ", 2)

# ---------------------------------------------------------------------
# Change 6: merge "Multimodal Online Student Engagement Dataset" + " "
#           into a single run
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Multimodal Online Student Engagement Dataset ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Multimodal Online Student Engagement Dataset ", 2)
